$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-11-02 04:41:33"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-11-02 04:41:19"
$wsZhCn.Range("K2").Value = "2016-11-02 04:42:13"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-11-02 04:41:33"
$wsDeDe.Range("K2").Value = "2016-11-02 04:42:31"
